# Update cryptocurrency price (D) and volume-change (E) columns
# to reflect the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.159.62"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.854.80"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.43"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4691"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2881"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06551"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.83"
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07968"
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.51"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.855.03"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.102"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6767"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.38"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("D17").Value = "30.140.97"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.59"
$ws.Range("E18").Value = "  +7.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007642"
$ws.Range("E19").Value = "  +4.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "2.097.22"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.201"
$ws.Range("E23").Value = "  -5.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.134"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.01"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.162"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.92"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.933"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.378"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09856"
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.464"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.295"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.022"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04692"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6981"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.703"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01865"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.322"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.32"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9982"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8380"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.25"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4132"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.205"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.015"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.88"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05652"
$ws.Range("E51").Value = "  +0.31%  "
